$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.708.10"
$ws.Range("E2").Value = "'  -6.81%  "

$ws.Range("D3").Value = "'3.686.59"
$ws.Range("E3").Value = "'  -6.23%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'576.96"
$ws.Range("E5").Value = "'  -4.02%  "

$ws.Range("D6").Value = "'171.28"
$ws.Range("E6").Value = "'  -0.29%  "

$ws.Range("D7").Value = "'3.680.57"
$ws.Range("E7").Value = "'  -6.18%  "

$ws.Range("D8").Value = "'0.622"
$ws.Range("E8").Value = "'  -8.93%  "

$ws.Range("E9").Value = "'  +0.06%  "

$ws.Range("D10").Value = "'0.701"
$ws.Range("E10").Value = "'  -10.68%  "

$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "'  -13.62%  "

$ws.Range("D12").Value = "'50.90"
$ws.Range("E12").Value = "'  -8.93%  "

$ws.Range("D13").Value = "'0.0000286"
$ws.Range("E13").Value = "'  -13.68%  "

$ws.Range("D14").Value = "'10.33"
$ws.Range("E14").Value = "'  -10.25%  "

$ws.Range("D15").Value = "'4.272.15"
$ws.Range("E15").Value = "'  -6.20%  "

$ws.Range("D16").Value = "'3.682.47"
$ws.Range("E16").Value = "'  -5.96%  "

$ws.Range("D17").Value = "'19.27"
$ws.Range("E17").Value = "'  -10.61%  "

$ws.Range("E18").Value = "'  -3.32%  "

$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "'  -9.84%  "

$ws.Range("E20").Value = "'  -9.70%  "

$ws.Range("D21").Value = "'67.417.16"
$ws.Range("E21").Value = "'  -7.05%  "

$ws.Range("D22").Value = "'403.18"
$ws.Range("E22").Value = "'  -8.86%  "

$ws.Range("D23").Value = "'4.46"
$ws.Range("E23").Value = "'  -5.77%  "

$ws.Range("D24").Value = "'87.58"
$ws.Range("E24").Value = "'  -8.15%  "

$ws.Range("D25").Value = "'3.02"
$ws.Range("E25").Value = "'  -9.66%  "

$ws.Range("D26").Value = "'12.64"
$ws.Range("E26").Value = "'  -9.94%  "

$ws.Range("D27").Value = "'10.74"
$ws.Range("E27").Value = "'  -3.47%  "

$ws.Range("D28").Value = "'6.03"
$ws.Range("E28").Value = "'  +1.82%  "

$ws.Range("D29").Value = "'3.77"
$ws.Range("E29").Value = "'  -12.74%  "

$ws.Range("D30").Value = "'9.40"
$ws.Range("E30").Value = "'  -9.63%  "

$ws.Range("D31").Value = "'32.39"
$ws.Range("E31").Value = "'  -9.44%  "

$ws.Range("D32").Value = "'7.37"
$ws.Range("E32").Value = "'  -6.66%  "

$ws.Range("D33").Value = "'12.33"
$ws.Range("E33").Value = "'  -10.93%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'605.70"
$ws.Range("E34").Value = "'  -4.42%  "

$ws.Range("D35").Value = "'64.68"
$ws.Range("E35").Value = "'  -5.78%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.115"
$ws.Range("E36").Value = "'  -9.81%  "

$ws.Range("D37").Value = "'42.80"
$ws.Range("E37").Value = "'  -15.47%  "

$ws.Range("D38").Value = "'0.0₃0879"
$ws.Range("E38").Value = "'  -11.98%  "

$ws.Range("E39").Value = "'  +0.12%  "

$ws.Range("D40").Value = "'0.391"
$ws.Range("E40").Value = "'  -8.11%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "'  -0.01%  "

$ws.Range("E42").Value = "'  -7.83%  "

$ws.Range("D43").Value = "'2.75"
$ws.Range("E43").Value = "'  +5.40%  "

$ws.Range("D44").Value = "'2.96"
$ws.Range("E44").Value = "'  -12.51%  "

$ws.Range("D45").Value = "'0.0431"
$ws.Range("E45").Value = "'  -9.69%  "

$ws.Range("D46").Value = "'2.85"
$ws.Range("E46").Value = "'  -12.86%  "

$ws.Range("D47").Value = "'9.17"
$ws.Range("E47").Value = "'  -12.85%  "

$ws.Range("D48").Value = "'2.809.98"
$ws.Range("E48").Value = "'  -0.41%  "

$ws.Range("E49").Value = "'  -6.20%  "

$ws.Range("E50").Value = "'  -9.73%  "

$ws.Range("D51").Value = "'3.12"
$ws.Range("E51").Value = "'  -7.15%  "
